# Adds the "Numero.Comorbidades (%)" block of rows to the comorbidities
# table, right after the "n" / "364" row and before "Artrite.Reumatoide".
#
# Structure per new row (matches the existing table rows):
#   cell 1: pStyle "Compact", centered, bold label run
#   cell 2: pStyle "Compact", centered, plain value run
# except for the header row "Numero.Comorbidades (%)" whose second cell
# stays an empty paragraph (pStyle "Compact" only, no run).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$w_ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newRows = @(
    @{ Label = "Numero.Comorbidades (%)"; Value = $null },
    @{ Label = "0"; Value = "73 (20.1)" },
    @{ Label = "1"; Value = "160 (44.0)" },
    @{ Label = "2"; Value = "105 (28.8)" },
    @{ Label = "3"; Value = "21 ( 5.8)" },
    @{ Label = "4"; Value = "4 ( 1.1)" },
    @{ Label = "6"; Value = "1 ( 0.3)" }
)

# New rows are inserted one by one right after the "n" row (table row 2),
# each insertion pushing the following rows further down, so they end up
# in the table in the same order as $newRows.
$insertAt = 3

foreach ($row in $newRows) {
    $t.Rows.Add($t.Rows.Item($insertAt)) | Out-Null

    # A freshly-added row starts life as a single merged cell spanning the
    # whole row width; split it into the table's real 2 columns.
    $t.Cell($insertAt, 1).Split(1, 2) | Out-Null

    $labelXml = '<w:p ' + $w_ns + '><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">' + $row.Label + '</w:t></w:r></w:p>'
    $t.Cell($insertAt, 1).Range.InsertXML($labelXml)

    if ($null -ne $row.Value) {
        $valueXml = '<w:p ' + $w_ns + '><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="center"/></w:pPr><w:r><w:t xml:space="preserve">' + $row.Value + '</w:t></w:r></w:p>'
    } else {
        $valueXml = '<w:p ' + $w_ns + '><w:pStyle w:val="Compact"/></w:p>'
    }
    $t.Cell($insertAt, 2).Range.InsertXML($valueXml)

    $insertAt = $insertAt + 1
}
